# Messwerte v5 - add derived "ei AD" / "ei Phil" / "Quanfehler" columns and
# the v5 samplerate block, per the commit "python fertig aber muss zugeordnet werden"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 2: new columns F (ei AD), G (ei Phil), I (Quanfehler) ---
$ws.Range("F2").Value = "ei AD"
$ws.Range("G2").Value = "ei Phil"
$ws.Range("I2").Value = "Quanfehler"

# --- Row 3 (first data row): standalone formulas (not part of the shared group below) ---
$ws.Range("F3").Formula = "=E3-C3"
$ws.Range("G3").Formula = "=E3-D3"
$ws.Range("I3").Formula = "=20/2^11"

# --- Rows 4-12: shared formulas ---
$ws.Range("F4:F12").Formula = "=E4-C4"
$ws.Range("G4:G12").Formula = "=E4-D4"

# --- Row 13: standard deviation summary for F/G ---
$ws.Range("E13").Value = "standardabweichung"
$ws.Range("F13").Formula = "=SQRT(SUMSQ(F3:F12)/9)"
$ws.Range("G13").Formula = "=SQRT(SUMSQ(G3:G12)/9)"

# --- Number formats for the new F/G columns ---
$ws.Range("F3:F13").NumberFormat = "0.00000000"
$ws.Range("G13").NumberFormat = "0.00000000"

# --- Row 14: label for the Quanfehler block next to v3/Oszil ---
$ws.Range("I14").Value = "Quanfehler"

# --- Row 15 (first row of second table): standalone formula + Quanfehler formula ---
$ws.Range("E15").Formula = "=C15-D15"
$ws.Range("I15").Formula = "=5/2^10"

# --- Rows 16-24: shared formula ---
$ws.Range("E16:E24").Formula = "=C16-D16"

# --- Row 25: standard deviation summary for E ---
$ws.Range("E25").Formula = "=SQRT(SUMSQ(E15:E24)/9)"
$ws.Range("E25").NumberFormat = "0.00000000"

# --- Number format for the E15:E24 differences ---
$ws.Range("E15:E24").NumberFormat = "0.000"

# --- v5 / Samplerate block ---
$ws.Range("A28").Value = "v5"
$ws.Range("B28").Value = "Samplerate"

# B29 must be stored as literal text "8021.0" (not converted to the number 8021),
# without leaving behind any extra number-format/style. Enter it as a text formula
# and then freeze it to a plain value via copy / paste-values.
$ws.Range("B29").Formula = '="8021.0"'
$ws.Range("B29").Copy()
$ws.Range("B29").PasteSpecial(-4163)

# --- Selection / view matches the saved state (also drops the stale topLeftCell) ---
$ws.Range("I14").Select()
